$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.680.64'
$ws.Range('E2').Value = '  -0.10%  '
$ws.Range('D3').Value = '1.918.44'
$ws.Range('E3').Value = '  +1.39%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').Value = "'239.97"
$ws.Range('E5').Value = '  -2.12%  '
$ws.Range('E6').Value = '  +0.08%  '
$ws.Range('D7').Value = "'0.4938"
$ws.Range('E7').Value = '  +0.23%  '
$ws.Range('D8').Value = "'0.2983"
$ws.Range('E8').Value = '  +0.88%  '
$ws.Range('D9').Value = "'0.06772"
$ws.Range('E9').Value = '  -0.41%  '
$ws.Range('D10').Value = '1.907.61'
$ws.Range('E10').Value = '  +0.81%  '
$ws.Range('D11').Value = "'17.18"
$ws.Range('E11').Value = '  -0.31%  '
$ws.Range('D12').Value = "'0.07351"
$ws.Range('E12').Value = '  +1.44%  '
$ws.Range('D13').Value = "'5.170"
$ws.Range('E13').Value = '  +2.43%  '
$ws.Range('D14').Value = "'88.91"
$ws.Range('E14').Value = '  -2.26%  '
$ws.Range('D15').Value = "'0.6726"
$ws.Range('E15').Value = '  -0.92%  '
$ws.Range('D16').Value = '30.657.07'
$ws.Range('E16').Value = '  -0.09%  '
$ws.Range('D17').Value = "'0.000007953"
$ws.Range('E17').Value = '  -0.41%  '
$ws.Range('D18').Value = "'13.54"
$ws.Range('E18').Value = '  +2.84%  '
$ws.Range('D19').Value = "'1.002"
$ws.Range('E19').Value = '  +0.13%  '
$ws.Range('D20').Value = '2.156.32'
$ws.Range('E20').Value = '  +1.16%  '
$ws.Range('D21').Value = "'5.366"
$ws.Range('E21').Value = '  +11.26%  '
$ws.Range('E22').Value = '  +0.15%  '
$ws.Range('D23').Value = "'201.86"
$ws.Range('E23').Value = '  +6.17%  '
$ws.Range('D24').Value = "'6.313"
$ws.Range('E24').Value = '  +2.72%  '
$ws.Range('D25').Value = "'9.648"
$ws.Range('E25').Value = '  +2.51%  '
$ws.Range('D26').Value = "'165.64"
$ws.Range('E26').Value = '  +6.31%  '
$ws.Range('D27').Value = "'18.89"
$ws.Range('E27').Value = '  -1.14%  '
$ws.Range('D28').Value = "'1.964"
$ws.Range('E28').Value = '  +3.34%  '
$ws.Range('E29').Value = '  +5.82%  '
$ws.Range('D30').Value = "'4.381"
$ws.Range('E30').Value = '  +0.83%  '
$ws.Range('D31').Value = "'0.09201"
$ws.Range('E31').Value = '  +1.28%  '
$ws.Range('D32').Value = "'4.069"
$ws.Range('E32').Value = '  +1.33%  '
$ws.Range('D33').Value = "'0.05285"
$ws.Range('E33').Value = '  +1.43%  '
$ws.Range('D34').Value = "'0.7428"
$ws.Range('E34').Value = '  -1.18%  '
$ws.Range('D35').Value = "'1.119"
$ws.Range('E35').Value = '  +0.89%  '
$ws.Range('D36').Value = "'2.730"
$ws.Range('E36').Value = '  -1.71%  '
$ws.Range('D37').Value = "'0.01841"
$ws.Range('E37').Value = '  +0.16%  '
$ws.Range('D38').Value = "'2.720"
$ws.Range('E38').Value = '  +1.31%  '
$ws.Range('D39').Value = "'0.9262"
$ws.Range('E39').Value = '  -1.26%  '
$ws.Range('D40').Value = "'2.080"
$ws.Range('E40').Value = '  -3.03%  '
$ws.Range('B41').Value = 'Aave'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D41').Value = "'73.47"
$ws.Range('E41').Value = '  +27.47%  '
$ws.Range('B42').Value = 'TheSandbox'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D42').Value = "'0.4461"
$ws.Range('E42').Value = '  +0.86%  '
$ws.Range('D43').Value = "'5.995"
$ws.Range('E43').Value = '  +3.98%  '
$ws.Range('D44').Value = "'106.77"
$ws.Range('E44').Value = '  +1.35%  '
$ws.Range('E45').Value = '  +0.13%  '
$ws.Range('D46').Value = "'0.1388"
$ws.Range('E46').Value = '  +3.44%  '
$ws.Range('D47').Value = "'7.643"
$ws.Range('E47').Value = '  +0.70%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = "'9.096"
$ws.Range('E48').Value = '  +3.90%  '
$ws.Range('B49').Value = 'Elrond'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range('D49').Value = "'35.48"
$ws.Range('E49').Value = '  +5.57%  '
$ws.Range('D50').Value = "'0.05881"
$ws.Range('E50').Value = '  +0.22%  '
$ws.Range('D51').Value = "'0.4044"
$ws.Range('E51').Value = '  +2.75%  '
